$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values to reflect repulled data / mean calculation
$ws.Range("F2").Value = -1
$ws.Range("F4").Value = -9
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -5
